$d = $word.ActiveDocument

# Locate the paragraph that ends "...LA razón Wo es {{catwosubrt}} con {{valorwosubrt}} vs LW {{valorwosubrtlw}}"
# so the new "level 3" paragraph can be inserted right after it (and before "En Partner cerramos...").
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("valorwosubrtlw}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $findRange.Paragraphs.Item(1).Index

$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a brand-new paragraph right after the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter() | Out-Null

# The freshly created paragraph is now the next one after the anchor; fill in its text.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "Dentro de level 3 se identifico que la razón wo es {{catwortlv3}} con {{valorwortlv3}} vs LW {{valorwortlv3lw}}"

# Match the "w:spacing w:after=0" formatting used by every other paragraph in this template.
$newPara.Format.SpaceAfter = 0
